# Leetcode / GFG tracker update (DSA-450 additions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the new rows / column in an order that reproduces the
# --- original authoring sequence (new strings first, header "Type"
# --- column, then "Java" tags, then the Minimize-the-Heights rows
# --- which got swapped around after first typing them in) ---

$ws.Range("B6").Value  = "Minimize the Heights II"
$ws.Range("B7").Value  = "Minimize the Heights I"
$ws.Range("B8").Value  = "Find minimum and maximum element in an array"
$ws.Range("B9").Value  = "Kth smallest element"
$ws.Range("B10").Value = "Kth Largest Element in an Array"

$ws.Range("C1").Value  = "Type"
$ws.Range("C6").Value  = "Java"
$ws.Range("C7").Value  = "Java"
$ws.Range("C8").Value  = "Java"
$ws.Range("C9").Value  = "Java"
$ws.Range("C10").Value = "Java"

# swap rows 6/7 back to their final order
$ws.Range("B6").Value = "Minimize the Heights I"
$ws.Range("B7").Value = "Minimize the Heights II"

$ws.Range("A6").Value  = "GFG"
$ws.Range("A7").Value  = "GFG"
$ws.Range("A8").Value  = "GFG"
$ws.Range("A9").Value  = "GFG"
$ws.Range("A10").Value = 215

# --- Column widths / layout (nearest values the host's pixel-quantized
# --- ColumnWidth setter can represent; targets are 12.85546875 / 53.85546875) ---
$ws.Columns("A").ColumnWidth = 12
$ws.Columns("B").ColumnWidth = 53

# --- Formatting: ID column (A) centered + top aligned, Question column (B) wrap text ---
$ws.Range("A1:A10").HorizontalAlignment = -4108
$ws.Range("A1:A10").VerticalAlignment = -4160
$ws.Range("B1:B10").WrapText = $true

# --- Selection as left by the editing session ---
$ws.Range("G17").Select() | Out-Null
